$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 101, copying row 101's formatting (style, height)
# so the new rows inherit style "15" and the 32.25pt custom row height, matching
# the rest of the data table.
$ws.Rows.Item(101).Copy() | Out-Null
$ws.Rows.Item(102).Insert(-4121) | Out-Null   # xlShiftDown
$ws.Rows.Item(101).Copy() | Out-Null
$ws.Rows.Item(103).Insert(-4121) | Out-Null   # xlShiftDown
$excel.CutCopyMode = $false

$fis = @"
/**
 * A <code>FileInputStream</code> obtains input bytes
 * from a file in a file system. What files
 * are  available depends on the host environment.
 *
 * <p><code>FileInputStream</code> is meant for reading streams of raw bytes
 * such as image data. For reading streams of characters, consider using
 * <code>FileReader</code>.
 *
 * @author  Arthur van Hoff
 * @see     java.io.File
 * @see     java.io.FileDescriptor
 * @see     java.io.FileOutputStream
 * @see     java.nio.file.Files#newInputStream
 * @since   JDK1.0
 */
public class FileInputStream extends InputStream ...
"@

$skip = @"
    /**
     * Skips over and discards <code>n</code> bytes of data from the
     * input stream.
     *
     * <p>The <code>skip</code> method may, for a variety of
     * reasons, end up skipping over some smaller number of bytes,
     * possibly <code>0</code>. If <code>n</code> is negative, an
     * <code>IOException</code> is thrown, even though the <code>skip</code>
     * method of the {@link InputStream} superclass does nothing in this case.
     * The actual number of bytes skipped is returned.
     *
     * <p>This method may skip more bytes than are remaining in the backing
     * file. This produces no exception and the number of bytes skipped
     * may include some number of bytes that were beyond the EOF of the
     * backing file. Attempting to read from the stream after skipping past
     * the end will result in -1 indicating the end of the file.
     *
     * @param      n   the number of bytes to be skipped.
     * @return     the actual number of bytes skipped.
     * @exception  IOException  if n is negative, if the stream does not
     *             support seek, or if an I/O error occurs.
     */
    public native long skip(long n) throws IOException...
"@

# Row 102: "comment style (class)" example -> FileInputStream class javadoc
# Row 103: "comment style (method)" example -> InputStream#skip method javadoc
$ws.Cells.Item(102, 1).Value = "java"
$ws.Cells.Item(103, 1).Value = "java"
$ws.Cells.Item(102, 2).Value = "comment style (class)"
$ws.Cells.Item(103, 2).Value = "comment style (method)"
$ws.Cells.Item(103, 3).Value = $skip
$ws.Cells.Item(102, 3).Value = $fis

# Re-assert the table's standard row height (the long javadoc text would
# otherwise auto-fit the rows taller than the rest of the table).
$ws.Rows.Item(102).RowHeight = 32.25
$ws.Rows.Item(103).RowHeight = 32.25

# Keep the sheet view in sync with where Excel would land the selection
# after typing the last new entry (one cell below/right of the new block).
$ws.Range("B106").Select() | Out-Null
